# Corrected values for T1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.58333
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 45385.3447201341
$ws.Range("G2").Value = 45387.3447201341

$ws.Range("E3").Value = 1.5
$ws.Range("F3").Value = 45387.3447201341
$ws.Range("G3").Value = 45389.3447201341
